$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$styleRef = $ws.Range("B2")  # known default-style cell, used to normalize styles on text-like-number cells

$ws.Range("D2").Value = "27.867.76"
$ws.Range("E2").Value = "  +3.24%  "
$ws.Range("D3").Value = "1.726.83"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "'217.01"
$ws.Range("D5").Style = $styleRef.Style
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'0.523"
$ws.Range("D6").Style = $styleRef.Style
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = $styleRef.Style
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "'23.98"
$ws.Range("D8").Style = $styleRef.Style
$ws.Range("E8").Value = "  +12.05%  "
$ws.Range("E9").Value = "  +4.59%  "
$ws.Range("D10").Value = "'0.0631"
$ws.Range("D10").Style = $styleRef.Style
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").Value = "'0.0901"
$ws.Range("D11").Style = $styleRef.Style
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Value = "1.966.81"
$ws.Range("D13").Value = "1.725.74"
$ws.Range("E13").Value = "  +3.69%  "
$ws.Range("D14").Value = "'4.24"
$ws.Range("D14").Style = $styleRef.Style
$ws.Range("E14").Value = "  +3.33%  "
$ws.Range("D15").Value = "'0.566"
$ws.Range("D15").Style = $styleRef.Style
$ws.Range("E15").Value = "  +6.14%  "
$ws.Range("D16").Value = "'68.10"
$ws.Range("D16").Style = $styleRef.Style
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").Value = "27.853.40"
$ws.Range("E17").Value = "  +3.18%  "
$ws.Range("D18").Value = "'241.65"
$ws.Range("D18").Style = $styleRef.Style
$ws.Range("E18").Value = "  +2.93%  "
$ws.Range("D19").Value = "'8.06"
$ws.Range("D19").Style = $styleRef.Style
$ws.Range("E19").Value = "  -1.36%  "
$ws.Range("D20").Value = "0.0₃0751"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("D21").Style = $styleRef.Style
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  +4.01%  "
$ws.Range("D23").Value = "'9.69"
$ws.Range("D23").Style = $styleRef.Style
$ws.Range("E23").Value = "  +4.82%  "
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").Value = "'148.64"
$ws.Range("D25").Style = $styleRef.Style
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").Value = "'7.55"
$ws.Range("D26").Style = $styleRef.Style
$ws.Range("E26").Value = "  +4.29%  "
$ws.Range("D27").Value = "'16.66"
$ws.Range("D27").Style = $styleRef.Style
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = $styleRef.Style
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "'0.0506"
$ws.Range("D30").Style = $styleRef.Style
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").Value = "'1.19"
$ws.Range("D31").Style = $styleRef.Style
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("D32").Value = "'3.45"
$ws.Range("D32").Style = $styleRef.Style
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("E33").Value = "  +4.41%  "
$ws.Range("D34").Value = "1.478.87"
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").Value = "'0.967"
$ws.Range("D36").Style = $styleRef.Style
$ws.Range("E36").Value = "  +6.34%  "
$ws.Range("D37").Value = "'0.615"
$ws.Range("D37").Style = $styleRef.Style
$ws.Range("E37").Value = "  +3.99%  "
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("E40").Value = "  +2.88%  "
$ws.Range("D41").Value = "'71.97"
$ws.Range("D41").Style = $styleRef.Style
$ws.Range("E41").Value = "  +6.57%  "
$ws.Range("D42").Value = "'5.87"
$ws.Range("D42").Style = $styleRef.Style
$ws.Range("E42").Value = "  +6.23%  "
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "1.871.05"
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("D46").Value = "'0.789"
$ws.Range("D46").Style = $styleRef.Style
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "'1.69"
$ws.Range("D47").Style = $styleRef.Style
$ws.Range("E47").Value = "  +10.04%  "
$ws.Range("D48").Value = "'91.75"
$ws.Range("D48").Style = $styleRef.Style
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("D49").Value = "0.0₆0110"
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.32"
$ws.Range("D50").Style = $styleRef.Style
$ws.Range("E50").Value = "  +3.86%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.106"
$ws.Range("D51").Style = $styleRef.Style
$ws.Range("E51").Value = "  +2.56%  "
